# Auto-generated Excel COM-interop script
# Applies numeric updates to columns H-N across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 156.59091
$ws.Range("I33").Value = 122.5
$ws.Range("J33").Value = 310
$ws.Range("K33").Value = 122.5
$ws.Range("L33").Value = 310
$ws.Range("M33").Value = 106.5
$ws.Range("N33").Value = -768
$ws.Range("H113").Value = 3290.5557
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 3290.5557
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 3290.5557
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -9798.555700000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 953.2083
$ws.Range("I2").Value = 871.6818
$ws.Range("J2").Value = 1850
$ws.Range("K2").Value = 871.6818
$ws.Range("L2").Value = 1850
$ws.Range("M2").Value = -758.6818
$ws.Range("N2").Value = -2076
$ws.Range("H61").Value = 2153.65
$ws.Range("I61").Value = 989.9286
$ws.Range("J61").Value = 4869
$ws.Range("K61").Value = 989.9286
$ws.Range("L61").Value = 4869
$ws.Range("M61").Value = -777.9286
$ws.Range("N61").Value = -5293
$ws.Range("H62").Value = 27924.375
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 27924.375
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 27924.375
$ws.Range("N62").Value = -29172.375
$ws.Range("H65").Value = 27924.375
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 27924.375
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 83773.125
$ws.Range("N65").Value = -90013.125
$ws.Range("H102").Value = 1360
$ws.Range("I102").Value = 1420
$ws.Range("J102").Value = 1000
$ws.Range("K102").Value = 1420
$ws.Range("L102").Value = 1000
$ws.Range("M102").Value = 202
$ws.Range("N102").Value = -4244
$ws.Range("H110").Value = 1152.8334
$ws.Range("I110").Value = 983.4
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 983.4
$ws.Range("L110").Value = 2000
$ws.Range("M110").Value = 1061.6
$ws.Range("H116").Value = 953.2083
$ws.Range("I116").Value = 871.6818
$ws.Range("J116").Value = 1850
$ws.Range("K116").Value = 871.6818
$ws.Range("L116").Value = 1850
$ws.Range("M116").Value = 1422.3182
$ws.Range("N116").Value = -6438
$ws.Range("H136").Value = 2153.65
$ws.Range("I136").Value = 989.9286
$ws.Range("J136").Value = 4869
$ws.Range("K136").Value = 2969.7858
$ws.Range("L136").Value = 14607
$ws.Range("M136").Value = -419.7857999999997
$ws.Range("N136").Value = -19707

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 953.2083
$ws.Range("I3").Value = 871.6818
$ws.Range("J3").Value = 1850
$ws.Range("K3").Value = 871.6818
$ws.Range("L3").Value = 1850
$ws.Range("M3").Value = -757.6818
$ws.Range("N3").Value = -2078
$ws.Range("H105").Value = 1501.64
$ws.Range("I105").Value = 1496.8182
$ws.Range("J105").Value = 1537
$ws.Range("K105").Value = 1496.8182
$ws.Range("L105").Value = 1537
$ws.Range("M105").Value = 250.1818000000001
$ws.Range("N105").Value = -5031

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1854.1666
$ws.Range("I58").Value = 2040.6111
$ws.Range("J58").Value = 1574.5
$ws.Range("K58").Value = 2040.6111
$ws.Range("L58").Value = 1574.5
$ws.Range("M58").Value = -1837.6111
$ws.Range("N58").Value = -1980.5
$ws.Range("H105").Value = 633
$ws.Range("I105").Value = 633
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 633
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 1114
$ws.Range("H122").Value = 4102.1763
$ws.Range("I122").Value = 4070.7778
$ws.Range("J122").Value = 4137.5
$ws.Range("K122").Value = 12212.3334
$ws.Range("L122").Value = 12412.5
$ws.Range("M122").Value = -9762.3334
$ws.Range("N122").Value = -17312.5
$ws.Range("H132").Value = 1695.9706
$ws.Range("I132").Value = 1358.7693
$ws.Range("J132").Value = 2791.875
$ws.Range("K132").Value = 4076.3079
$ws.Range("L132").Value = 8375.625
$ws.Range("M132").Value = -1546.3079
$ws.Range("N132").Value = -13435.625
$ws.Range("H134").Value = 1045.7949
$ws.Range("I134").Value = 1032.697
$ws.Range("J134").Value = 1117.8334
$ws.Range("K134").Value = 3098.090999999999
$ws.Range("L134").Value = 3353.5002
$ws.Range("M134").Value = -563.0909999999994
$ws.Range("N134").Value = -8423.5002
$ws.Range("H136").Value = 1854.1666
$ws.Range("I136").Value = 2040.6111
$ws.Range("J136").Value = 1574.5
$ws.Range("K136").Value = 6121.8333
$ws.Range("L136").Value = 4723.5
$ws.Range("M136").Value = -3571.8333
$ws.Range("N136").Value = -9823.5
$ws.Range("H140").Value = 44950
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 44950
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 44950
$ws.Range("N140").Value = -55310

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 702
$ws.Range("I92").Value = 502
$ws.Range("J92").Value = 968.6667
$ws.Range("K92").Value = 1506
$ws.Range("L92").Value = 2906.0001
$ws.Range("M92").Value = -258
$ws.Range("N92").Value = -5402.0001
$ws.Range("H97").Value = 180
$ws.Range("I97").Value = 180
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 540
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -44
$ws.Range("N97").ClearContents()
$ws.Range("H131").Value = 935.7959
$ws.Range("I131").Value = 299.5
$ws.Range("J131").Value = 962.8723
$ws.Range("K131").Value = 898.5
$ws.Range("L131").Value = 2888.6169
$ws.Range("M131").Value = 4141.5
$ws.Range("N131").Value = -12968.6169

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1225
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1225
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1225
$ws.Range("N113").Value = -5565
$ws.Range("H122").Value = 2963.9412
$ws.Range("I122").Value = 1050.875
$ws.Range("J122").Value = 4664.4443
$ws.Range("K122").Value = 3152.625
$ws.Range("L122").Value = 13993.3329
$ws.Range("M122").Value = -702.625
$ws.Range("N122").Value = -18893.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2353.25
$ws.Range("I61").Value = 1265.1
$ws.Range("J61").Value = 2957.7778
$ws.Range("K61").Value = 1265.1
$ws.Range("L61").Value = 2957.7778
$ws.Range("M61").Value = -1063.1
$ws.Range("N61").Value = -3361.7778
$ws.Range("H113").Value = 2353.25
$ws.Range("I113").Value = 1265.1
$ws.Range("J113").Value = 2957.7778
$ws.Range("K113").Value = 1265.1
$ws.Range("L113").Value = 2957.7778
$ws.Range("M113").Value = 904.9000000000001
$ws.Range("N113").Value = -7297.7778
$ws.Range("H132").Value = 5020.156
$ws.Range("I132").Value = 2775.1428
$ws.Range("J132").Value = 8717.823
$ws.Range("K132").Value = 8325.428400000001
$ws.Range("L132").Value = 26153.469
$ws.Range("M132").Value = -5795.428400000001
$ws.Range("N132").Value = -31213.469
$ws.Range("H136").Value = 3207.795
$ws.Range("I136").Value = 2916.16
$ws.Range("J136").Value = 3728.5715
$ws.Range("K136").Value = 8748.48
$ws.Range("L136").Value = 11185.7145
$ws.Range("M136").Value = -6198.48
$ws.Range("N136").Value = -16285.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2626.7
$ws.Range("I62").Value = 1973.6
$ws.Range("J62").Value = 3279.8
$ws.Range("K62").Value = 1973.6
$ws.Range("L62").Value = 3279.8
$ws.Range("M62").Value = -1349.6
$ws.Range("N62").Value = -4527.8
$ws.Range("H65").Value = 2626.7
$ws.Range("I65").Value = 1973.6
$ws.Range("J65").Value = 3279.8
$ws.Range("K65").Value = 9868
$ws.Range("L65").Value = 16399
$ws.Range("M65").Value = -6748
$ws.Range("N65").Value = -22639
$ws.Range("H132").Value = 967.06976
$ws.Range("I132").Value = 637.7308
$ws.Range("J132").Value = 1470.7646
$ws.Range("K132").Value = 1913.1924
$ws.Range("L132").Value = 4412.293799999999
$ws.Range("M132").Value = 616.8075999999999
$ws.Range("N132").Value = -9472.293799999999
$ws.Range("H136").Value = 868.25806
$ws.Range("I136").Value = 685.0526
$ws.Range("J136").Value = 1158.3334
$ws.Range("K136").Value = 2055.1578
$ws.Range("L136").Value = 3475.0002
$ws.Range("M136").Value = 494.8422
$ws.Range("N136").Value = -8575.0002

